$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new scene rows at the bottom of the table
$ws.Range("A21").Value = 13010017
$ws.Range("B21").Value = 1149
$ws.Range("C21").Value = 584
$ws.Range("D21").Value = 75
$ws.Range("E21").Value = 58
$ws.Range("F21").Value = "lcxj"

$ws.Range("A22").Value = 13010018
$ws.Range("B22").Value = 840
$ws.Range("C22").Value = 444
$ws.Range("D22").Value = 74
$ws.Range("E22").Value = 66
$ws.Range("F22").Value = "ygld"

# Grow the worksheet table so the new rows are included
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F22"))

# Give column A an explicit width like the other formatted columns
$ws.Columns.Item(1).ColumnWidth = 8.71

# Leave the selection where the author left it after typing the new data
$ws.Range("C22").Select() | Out-Null
